$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.041.49"
$ws.Range("E2").Value = "  -6.89%  "
$ws.Range("D3").Value = "2.549.96"
$ws.Range("E3").Value = "  -2.67%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.18%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "298.35"
$ws.Range("E5").Value = "  -4.04%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "94.23"
$ws.Range("E6").Value = "  -5.45%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.575"
$ws.Range("E7").Value = "  -3.81%  "
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.550"
$ws.Range("E9").Value = "  -5.82%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "36.08"
$ws.Range("E10").Value = "  -7.85%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0810"
$ws.Range("E11").Value = "  -4.42%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.74"
$ws.Range("E12").Value = "  -4.76%  "
$ws.Range("E13").Value = "  +1.22%  "
$ws.Range("D14").Value = "2.939.42"
$ws.Range("E14").Value = "  -2.72%  "
$ws.Range("D15").Value = "2.560.11"
$ws.Range("E15").Value = "  -2.22%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.872"
$ws.Range("E16").Value = "  -5.28%  "
$ws.Range("E17").Value = "  -4.88%  "
$ws.Range("D18").Value = "43.051.95"
$ws.Range("E18").Value = "  -7.42%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.96"
$ws.Range("E19").Value = "  +0.98%  "
$ws.Range("D20").Value = "0.0₃0980"
$ws.Range("E20").Value = "  -4.05%  "
$ws.Range("E21").Value = "  -2.08%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "71.77"
$ws.Range("E22").Value = "  -2.23%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "260.52"
$ws.Range("E23").Value = "  -11.42%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.92"
$ws.Range("E24").Value = "  -4.99%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "29.51"
$ws.Range("E25").Value = "  -1.24%  "
$ws.Range("E26").Value = "  -5.46%  "
$ws.Range("E27").Value = "  +0.00%  "
$ws.Range("E28").Value = "  -7.35%  "
$ws.Range("E29").Value = "  -4.10%  "
$ws.Range("E31").Value = "  -5.15%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "155.03"
$ws.Range("E32").Value = "  -2.30%  "
$ws.Range("E33").Value = "  -3.54%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.38"
$ws.Range("E34").Value = "  -6.19%  "
$ws.Range("E35").Value = "  -2.64%  "
$ws.Range("E36").Value = "  -5.36%  "
$ws.Range("E37").Value = "  -5.83%  "
$ws.Range("E38").Value = "  -3.32%  "

# Row 39/40: swap coin identities (EnergySwap <-> Celestia) with updated values
$ws.Range("B39").Value = "Celestia"
$ws.Range("C39").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "16.51"
$ws.Range("E39").Value = "  +4.46%  "
$ws.Range("B40").Value = "EnergySwap"
$ws.Range("C40").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "23.17"
$ws.Range("E40").Value = "  +6.96%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.47"
$ws.Range("E41").Value = "  -3.65%  "
$ws.Range("E42").Value = "  -5.44%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.88"
$ws.Range("E43").Value = "  -4.07%  "
$ws.Range("D44").Value = "2.068.05"
$ws.Range("E44").Value = "  -2.99%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.998"
$ws.Range("E45").Value = "  -0.11%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "85.71"
$ws.Range("E46").Value = "  -11.95%  "
$ws.Range("E47").Value = "  +2.79%  "
$ws.Range("D48").Value = "2.796.33"
$ws.Range("E48").Value = "  -2.73%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.75"
$ws.Range("E49").Value = "  -8.11%  "
$ws.Range("E50").Value = "  -2.75%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "104.24"

